$d = $word.ActiveDocument

# Update the change-of-plea / payment / suspension date from June 04, 2022 to June 05, 2022
$d.Content.Find.Execute("June 04, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "June 05, 2022", 2)

# Update the community service completion deadline from August 03, 2022 to August 04, 2022
$d.Content.Find.Execute("August 03, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "August 04, 2022", 2)
